# Exercise 4: rename the main data sheet and add YearlyIncome / MonthlyIncome
# summary-table tabs to the BikeStoreSample workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename the original (and only) sheet from "Sheet1" to "OrderDetailsData"
# ---------------------------------------------------------------------------
$wsData = $wb.Worksheets.Item(1)
$wsData.Name = "OrderDetailsData"

# Reset the view back to the top-left corner / A1 (no stale scroll position)
$wsData.Activate()
$wsData.Range("A1").Select() | Out-Null

# ---------------------------------------------------------------------------
# 2. Add the "YearlyIncome" sheet after OrderDetailsData
# ---------------------------------------------------------------------------
$wsYear = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsData)
$wsYear.Name = "YearlyIncome"

$wsYear.Range("A1").Value = "Year"
$wsYear.Range("B1").Value = 2016
$wsYear.Range("C1").Value = 2017
$wsYear.Range("D1").Value = 2018
$wsYear.Range("A2").Value = "Total "

$wsYear.Columns.Item(1).ColumnWidth = 10
$wsYear.Columns.Item(2).ColumnWidth = 13.57
$wsYear.Columns.Item(3).ColumnWidth = 17.43
$wsYear.Columns.Item(4).ColumnWidth = 17

$wsYear.Range("A1:D1").Font.Bold = $true
$wsYear.Range("A1:D2").Borders.LineStyle = 1
$wsYear.Range("B2:D2").NumberFormat = "_(""$""* #,##0.00_);_(""$""* \(#,##0.00\);_(""$""* ""-""??_);_(@_)"

# ---------------------------------------------------------------------------
# 3. Add the "MonthlyIncome" sheet after YearlyIncome
# ---------------------------------------------------------------------------
$wsMonth = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsYear)
$wsMonth.Name = "MonthlyIncome"

$wsMonth.Range("A1").Value = "Month"
$wsMonth.Range("B1").Value = 2016
$wsMonth.Range("C1").Value = 2017
$wsMonth.Range("D1").Value = 2018

for ($m = 1; $m -le 12; $m++) {
    $wsMonth.Cells.Item($m + 1, 1).Value = $m
}
$wsMonth.Range("A14").Value = "Total"

$wsMonth.Columns.Item(1).ColumnWidth = 10
$wsMonth.Columns.Item(2).ColumnWidth = 13.57
$wsMonth.Columns.Item(3).ColumnWidth = 17.43
$wsMonth.Columns.Item(4).ColumnWidth = 17

$wsMonth.Range("A1:D1").Font.Bold = $true
$wsMonth.Range("A14").Font.Bold = $true
$wsMonth.Range("A1:D14").Borders.LineStyle = 1
$wsMonth.Range("B2:D14").NumberFormat = "_(""$""* #,##0.00_);_(""$""* \(#,##0.00\);_(""$""* ""-""??_);_(@_)"

# ---------------------------------------------------------------------------
# 4. Leave OrderDetailsData as the active sheet/tab
# ---------------------------------------------------------------------------
$wsData.Activate()
